# Update countries & provincias Spain
# Refreshes the COVID-19 country ranking table ("Pais" sheet) with the
# 15-Jun-2020 00:05 data snapshot. The table is sorted descending by
# "Casos totales" (column B); a handful of countries (Costa de Marfil,
# Surinam, Groenlandia, Islas Turcas y Caicos) moved up in rank, which
# cascades their new figures into rows that used to hold the country
# immediately below them. Row numbers/positions on the sheet do not
# change - only the country name + stats shown in certain rows do.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp caption at the bottom of the sheet
# Row 1: Datos actualizados a 14 de Junio de 2020 a las 22:48 -> Datos actualizados a 15 de Junio de 2020 a las 00:05
$ws.Range("A1").Value = "Datos actualizados a 15 de Junio de 2020 a las 00:05"

# Row 4: Estados Unidos -> Estados Unidos
$ws.Range("B4").Value = 2160959
$ws.Range("C4").Value = 18735
$ws.Range("D4").Value = 860204
$ws.Range("E4").Value = 1182910
$ws.Range("G4").Value = 318
$ws.Range("H4").Value = 117845

# Row 5: Brasil -> Brasil
$ws.Range("B5").Value = 867624
$ws.Range("C5").Value = 16828
$ws.Range("E5").Value = 386780
$ws.Range("G5").Value = 541
$ws.Range("H5").Value = 43332

# Row 11: Peru -> Peru
$ws.Range("B11").Value = 229736
$ws.Range("C11").Value = 4604
$ws.Range("D11").Value = 115579
$ws.Range("E11").Value = 107469
$ws.Range("G11").Value = 190
$ws.Range("H11").Value = 6688

# Row 12: Alemania -> Alemania
$ws.Range("B12").Value = 187671
$ws.Range("C12").Value = 248
$ws.Range("E12").Value = 6601
$ws.Range("G12").Value = 3
$ws.Range("H12").Value = 8870

# Row 76: Uzbekistan -> Costa de Marfil
$ws.Range("A76").Value = "Costa de Marfil"
$ws.Range("B76").Value = 5084
$ws.Range("C76").Value = 236
$ws.Range("D76").Value = 2505
$ws.Range("E76").Value = 2534
$ws.Range("H76").Value = 45

# Row 77: Tayikistan -> Uzbekistan
$ws.Range("A77").Value = "Uzbekistan"
$ws.Range("B77").Value = 5080
$ws.Range("C77").Value = 114
$ws.Range("D77").Value = 3943
$ws.Range("E77").Value = 1118
$ws.Range("H77").Value = 19

# Row 78: Costa de Marfil -> Tayikistan
$ws.Range("A78").Value = "Tayikistan"
$ws.Range("B78").Value = 5035
$ws.Range("C78").Value = 64
$ws.Range("D78").Value = 3409
$ws.Range("E78").Value = 1576
$ws.Range("H78").Value = 50

# Row 90: Bulgaria -> Bulgaria
$ws.Range("B90").Value = 3290
$ws.Range("C90").Value = 24
$ws.Range("D90").Value = 1730
$ws.Range("E90").Value = 1386
$ws.Range("G90").Value = 2
$ws.Range("H90").Value = 174

# Row 159: Montenegro -> Montenegro
$ws.Range("B159").Value = 325
$ws.Range("C159").Value = 1
$ws.Range("E159").Value = 1

# Row 161: Martinica -> Surinam
$ws.Range("A161").Value = "Surinam"
$ws.Range("B161").Value = 208
$ws.Range("C161").Value = 12
$ws.Range("D161").Value = 9
$ws.Range("E161").Value = 196
$ws.Range("H161").Value = 3

# Row 162: Mongolia -> Martinica
$ws.Range("A162").Value = "Martinica"
$ws.Range("B162").Value = 202
$ws.Range("E162").Value = 90
$ws.Range("H162").Value = 14

# Row 163: Surinam -> Mongolia
$ws.Range("A163").Value = "Mongolia"
$ws.Range("B163").Value = 197
$ws.Range("D163").Value = 98
$ws.Range("E163").Value = 99
$ws.Range("H163").Value = 0

# Row 173: Angola -> Angola
$ws.Range("B173").Value = 140
$ws.Range("C173").Value = 2
$ws.Range("E173").Value = 73

# Row 206: Islas Malvinas -> Groenlandia
$ws.Range("A206").Value = "Groenlandia"

# Row 207: Groenlandia -> Islas Malvinas
$ws.Range("A207").Value = "Islas Malvinas"

# Row 208: Santa Sede -> Islas Turcas y Caicos
$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("D208").Value = 11
$ws.Range("H208").Value = 1

# Row 209: Islas Turcas y Caicos -> Santa Sede
$ws.Range("A209").Value = "Santa Sede"
$ws.Range("D209").Value = 12
$ws.Range("H209").Value = 0
